# Generate Report for Handoff
#
# Localization status refresh: the "Status" columns move from
# "Handed back: in sync with en-US" to "Ready for handoff", and the
# handoff timestamps tick forward a minute or so. The two status columns
# on the Overview sheet (and the now-shorter text) also get narrowed.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status
$wsZhCn.Range("C2").Value = "Ready for handoff"       # Status column
$wsDeDe.Range("C2").Value = "Ready for handoff"       # Status column

# --- Timestamps ---
$wsOverview.Range("G2").Value = "2016-08-25 00:56:17" # Latest HO Xliff Generate Date
$wsZhCn.Range("H2").Value = "2016-08-25 00:56:12"     # Latest Handoff Datetime (zh-cn)
$wsDeDe.Range("H2").Value = "2016-08-25 00:56:17"     # Latest Handoff Datetime (de-de), shares the Overview value

# --- Column widths: the Status columns shrink along with their new, shorter text ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333  # E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333  # F (de-de)
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333      # C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333      # C (Status)
